$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell reference -> new value, derived from the authoritative diff.
# Column D (Price) values are prefixed with a leading apostrophe so that
# Excel stores them as literal text (matching the original inlineStr
# cells) instead of silently reinterpreting them as numbers and dropping
# formatting such as trailing zeros ("1.00" -> 1) or thousands dots
# ("62.819.54").
$updates = @{
    "D2" = "'62.819.54"
    "E2" = "  +4.91%  "
    "D3" = "'3.473.07"
    "E3" = "  +4.13%  "
    "E4" = "  +0.01%  "
    "D5" = "'408.67"
    "E5" = "  -0.81%  "
    "D6" = "'131.16"
    "E6" = "  +17.85%  "
    "D7" = "'3.466.45"
    "E7" = "  +4.30%  "
    "D8" = "'0.598"
    "E8" = "  +2.76%  "
    "D9" = "'1.00"
    "E9" = "  +0.07%  "
    "D10" = "'0.692"
    "E10" = "  +9.99%  "
    "D11" = "'0.128"
    "E11" = "  +29.75%  "
    "D12" = "'42.78"
    "E12" = "  +7.75%  "
    "E13" = "  -0.85%  "
    "D14" = "'4.020.34"
    "E14" = "  +4.30%  "
    "D15" = "'8.72"
    "E15" = "  +2.74%  "
    "D16" = "'20.10"
    "E16" = "  +4.76%  "
    "D17" = "'3.465.63"
    "E17" = "  +4.68%  "
    "D18" = "'62.753.28"
    "E18" = "  +5.32%  "
    "D19" = "'1.05"
    "E19" = "  +0.50%  "
    "D20" = "'10.94"
    "E20" = "  +2.86%  "
    "D21" = "'0.0000136"
    "E21" = "  +23.49%  "
    "D22" = "'3.37"
    "E22" = "  +1.25%  "
    "D23" = "'82.54"
    "E23" = "  +9.69%  "
    "D24" = "'13.15"
    "E24" = "  +1.03%  "
    "D25" = "'309.41"
    "E25" = "  +2.35%  "
    "D26" = "'3.17"
    "E26" = "  -0.10%  "
    "D27" = "'30.36"
    "E27" = "  +6.47%  "
    "D28" = "'8.30"
    "E28" = "  +5.99%  "
    "D29" = "'7.76"
    "E29" = "  +4.57%  "
    "E30" = "  -1.55%  "
    "E31" = "  -2.10%  "
    "E32" = "  +4.49%  "
    "D33" = "'2.68"
    "E33" = "  +0.11%  "
    "D34" = "'11.94"
    "E34" = "  +3.59%  "
    "D35" = "'43.39"
    "E35" = "  +9.85%  "
    "E36" = "  +0.13%  "
    "D37" = "'0.0493"
    "E37" = "  -2.12%  "
    "D38" = "'52.59"
    "E38" = "  +1.65%  "
    "E39" = "  +5.29%  "
    "D40" = "'0.997"
    "E40" = "  -0.09%  "
    "D41" = "'3.01"
    "E41" = "  -4.71%  "
    "B42" = "Stellar"
    "C42" = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
    "D42" = "'0.126"
    "E42" = "  +2.71%  "
    "E43" = "  +3.97%  "
    "B44" = "Monero"
    "C44" = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
    "D44" = "'137.98"
    "E44" = "  -0.57%  "
    "D45" = "'17.55"
    "E45" = "  +4.66%  "
    "B46" = "NEARProtocol"
    "C46" = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
    "D46" = "'3.98"
    "E46" = "  +1.76%  "
    "B47" = "TheGraph"
    "C47" = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
    "D47" = "'0.286"
    "E47" = "  +0.76%  "
    "E48" = "  -0.70%  "
    "D49" = "'22.36"
    "E49" = "  +0.68%  "
    "D50" = "'2.207.97"
    "E50" = "  +0.58%  "
    "D51" = "'3.817.63"
    "E51" = "  +4.49%  "
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
